$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cross product of Sending cluster x Target cluster (ECs, FAPs, sCs) for Il16-Cd4,
# recomputed per Dr Hou's advice. Columns A:T per the existing header row.
$rows = @()
$rows += ,@("ECs","Il16","Cd4","ECs",3,1,7.933820999999999,23.801463,0.3933990084177585,0.3933990084177585,1,0.3333333333333333,1.635729666666667,4.907189,0.5314629201652572,0.5314629201652572,12.977586379723,116.798277417507,0.2090769858038185,0.2090769858038185)
$rows += ,@("ECs","Il16","Cd4","FAPs",3,1,7.933820999999999,23.801463,0.3933990084177585,0.3933990084177585,3,1,1.170126666666667,3.51038,0.3801844203860328,0.3801844203860328,9.28357552066,83.55217968593999,0.1495641739957456,0.1495641739957455)
$rows += ,@("ECs","Il16","Cd4","sCs",3,1,7.933820999999999,23.801463,0.3933990084177585,0.3933990084177585,2,0.6666666666666666,0.2719306666666667,0.8157920000000001,0.08835265944870997,0.08835265944870996,2.157449233744,19.417043103696,0.03475784861819441,0.0347578486181944)
$rows += ,@("FAPs","Il16","Cd4","ECs",3,1,11.58194133333333,34.745824,0.5742912823576415,0.5742912823576415,1,0.3333333333333333,1.635729666666667,4.907189,0.5314629201652572,0.5314629201652572,18.94492503652622,170.504325328736,0.3052145219472424,0.3052145219472424)
$rows += ,@("FAPs","Il16","Cd4","FAPs",3,1,11.58194133333333,34.745824,0.5742912823576415,0.5742912823576415,3,1,1.170126666666667,3.51038,0.3801844203860328,0.3801844203860328,13.55233840590222,121.97104565312,0.2183365983158914,0.2183365983158914)
$rows += ,@("FAPs","Il16","Cd4","sCs",3,1,11.58194133333333,34.745824,0.5742912823576415,0.5742912823576415,2,0.6666666666666666,0.2719306666666667,0.8157920000000001,0.08835265944870997,0.08835265944870996,3.149485028067556,28.345365252608,0.05074016209450764,0.05074016209450763)
$rows += ,@("sCs","Il16","Cd4","ECs",2,0.6666666666666666,0.6516016666666666,1.954805,0.03230970922460003,0.03230970922460003,1,0.3333333333333333,1.635729666666667,4.907189,0.5314629201652572,0.5314629201652572,1.065844177016111,9.592597593145,0.01717141241419628,0.01717141241419628)
$rows += ,@("sCs","Il16","Cd4","FAPs",2,0.6666666666666666,0.6516016666666666,1.954805,0.03230970922460003,0.03230970922460003,3,1,1.170126666666667,3.51038,0.3801844203860328,0.3801844203860328,0.7624564862111112,6.8621083759,0.01228364807439582,0.01228364807439582)
$rows += ,@("sCs","Il16","Cd4","sCs",2,0.6666666666666666,0.6516016666666666,1.954805,0.03230970922460003,0.03230970922460003,2,0.6666666666666666,0.2719306666666667,0.8157920000000001,0.08835265944870997,0.08835265944870996,0.1771904756177778,1.59471428056,0.00285464873600793,0.00285464873600793)

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
